$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number of data rows below header (rows 2..7 => 6 rows), 6 columns (A..F)
$numRows = 7
$numCols = 6

# Capture all original values (header row included) before overwriting anything,
# since the column permutation reads from multiple source columns.
$orig = @{}
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $orig["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# New column c (1-based) takes its values from old column map[c] (1-based)
# New A(1)=Old B(2), New B(2)=Old E(5), New C(3)=Old C(3),
# New D(4)=Old F(6), New E(5)=Old A(1), New F(6)=Old D(4)
$map = @{ 1 = 2; 2 = 5; 3 = 3; 4 = 6; 5 = 1; 6 = 4 }

for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $srcCol = $map[$c]
        $ws.Cells.Item($r, $c).Value = $orig["$r,$srcCol"]
    }
}
